# Weekly update: a new price observation is inserted as row 36
# ("Vega Modelo de Temuco", Maracuyá, fecha 2023-09-08), pushing the
# existing rows 36-119 down to 37-120 (dimension grows from T119 to T120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36, shifting rows 36:119 down to 37:120.
$ws.Rows("36").Insert()

# Populate the newly inserted row 36 with the new weekly observation.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = "2023-09-08"
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108003
$ws.Range("J36").Value = "Maracuyá"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 30
$ws.Range("N36").Value = 50000
$ws.Range("O36").Value = 50000
$ws.Range("P36").Value = 50000
$ws.Range("Q36").Value = "$/caja 18 kilos"
$ws.Range("R36").Value = "Región de Arica y Parinacota"
$ws.Range("S36").Value = 2778
$ws.Range("T36").Value = 18
